# Generate Report for Archive
#
# 1. Replace the shared "Ready for handoff" status text with "In Translation"
#    everywhere it is used (Overview!E2:F3, zh-cn!C2:C3, de-de!C2:C3).
# 2. Shrink the "Status" column widths (Overview columns E & F, and
#    column C on the zh-cn / de-de sheets) from ~17.22 to ~13.41 characters.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"
# ColumnWidth is quantized by the host to the nearest 1/6 of a character, so
# 12.5 is the input that lands on the stored width closest to 13.4101845877511.
$newWidth = 12.5

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        # Keep the string literal on the left of -eq so PowerShell compares
        # as strings instead of coercing (e.g. boolean-looking cell values
        # such as "True"/"False" would otherwise match any non-empty string).
        if ($oldStatus -eq $cell.Value2) {
            $cell.Value2 = $newStatus
        }
    }
}

$overview = $wb.Worksheets.Item("Overview")
$overview.Columns.Item(5).ColumnWidth = $newWidth
$overview.Columns.Item(6).ColumnWidth = $newWidth

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Columns.Item(3).ColumnWidth = $newWidth

$dede = $wb.Worksheets.Item("de-de")
$dede.Columns.Item(3).ColumnWidth = $newWidth
